$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Acad" semester label was a calendar-year error; the whole first-year
# block (rows 2-121 of the SemType column E) should read "Vacation".
# Rows 2-65 previously held "Acad" (shared string renamed to "Vacation"),
# rows 66-121 already held "Vacation".
$ws.Range("E2:E121").Value = "Vacation"

# The Exam flag (column H) for rows 38-57 was incorrectly marked 1 (exam)
# during that same vacation block; correct it back to 0.
$ws.Range("H38:H57").Value = 0
